$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Column widths (Excel stores width with a +0.8333 padding relative to
# the ColumnWidth property, so compensate to land on the exact target)
$ws.Columns.Item(4).ColumnWidth = 13 - 0.8333333333333334
$ws.Columns.Item(5).ColumnWidth = 22 - 0.8333333333333334

# Row 3 - 240X80 PORCELANATO
$ws.Cells.Item(3, 4).Value = 0
$ws.Cells.Item(3, 5).Value = 8668.91
$ws.Cells.Item(3, 6).Value = 0

# Row 4 - FREGADEROS DE COCINA
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = 372.993863046034
$ws.Cells.Item(4, 6).Value = 0

# Row 6 - GRIFERIAS
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(6, 5).Value = 106.82
$ws.Cells.Item(6, 6).Value = 0

# Row 7 - INODOROS
$ws.Cells.Item(7, 3).Value = 1600
$ws.Cells.Item(7, 4).Value = 0
$ws.Cells.Item(7, 5).Value = 1600
$ws.Cells.Item(7, 6).Value = 0

# Row 8 - LAVABOS
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = 625
$ws.Cells.Item(8, 6).Value = 0

# Row 13 - PANELES PU
$ws.Cells.Item(13, 3).Value = 130
$ws.Cells.Item(13, 5).Value = 130

# Row 14 - PANELES PVC
$ws.Cells.Item(14, 3).Value = 240
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = 240
$ws.Cells.Item(14, 6).Value = 0

# Row 16 - PORCELANATO
$ws.Cells.Item(16, 3).Value = 38756.54
$ws.Cells.Item(16, 4).Value = 6266.88
$ws.Cells.Item(16, 5).Value = 32489.66
$ws.Cells.Item(16, 6).Value = 0.1616986449254758

# Row 19 - TOTAL
$ws.Cells.Item(19, 3).Value = 58223.00386304603
$ws.Cells.Item(19, 4).Value = 6266.88
$ws.Cells.Item(19, 5).Value = 51956.12386304603
$ws.Cells.Item(19, 6).Value = 0.1076358068838418
